$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.517.55'
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").Value = '  +0.69%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.728.28'
$ws.Range("D3").Style = "Normal"

$ws.Range("E3").Value = '  +0.63%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '244.89'
$ws.Range("D5").Style = "Normal"

$ws.Range("E5").Value = '  +2.37%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("D6").Style = "Normal"

$ws.Range("E6").Value = '  +0.03%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4791'
$ws.Range("D7").Style = "Normal"

$ws.Range("E7").Value = '  +1.15%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2668'
$ws.Range("D8").Style = "Normal"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06223'
$ws.Range("D9").Style = "Normal"

$ws.Range("E9").Value = '  +0.23%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.730.73'
$ws.Range("D10").Style = "Normal"

$ws.Range("E10").Value = '  +0.77%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07159'
$ws.Range("D11").Style = "Normal"

$ws.Range("E11").Value = '  +1.44%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '15.67'
$ws.Range("D12").Style = "Normal"

$ws.Range("E12").Value = '  +2.31%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.6142'
$ws.Range("D13").Style = "Normal"

$ws.Range("E13").Value = '  +4.25%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.526'
$ws.Range("D14").Style = "Normal"

$ws.Range("E14").Value = '  +2.53%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '76.92'
$ws.Range("D15").Style = "Normal"

$ws.Range("E15").Value = '  +1.21%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.000'
$ws.Range("D16").Style = "Normal"

$ws.Range("E16").Value = '  +0.01%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.524.51'
$ws.Range("D17").Style = "Normal"

$ws.Range("E17").Value = '  +0.74%  '

$ws.Range("E18").Value = '  -0.01%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000006977'
$ws.Range("D19").Style = "Normal"

$ws.Range("E19").Value = '  +2.34%  '

$ws.Range("E20").Value = '  +0.89%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.952.36'
$ws.Range("D21").Style = "Normal"

$ws.Range("E21").Value = '  +0.85%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.525'
$ws.Range("D22").Style = "Normal"

$ws.Range("E22").Value = '  -0.25%  '

$ws.Range("E23").Value = '  +2.06%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.285'
$ws.Range("D24").Style = "Normal"

$ws.Range("E24").Value = '  -0.56%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '136.51'
$ws.Range("D25").Style = "Normal"

$ws.Range("E25").Value = '  +1.62%  '

$ws.Range("E26").Value = '  +0.65%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.791'
$ws.Range("D27").Style = "Normal"

$ws.Range("E27").Value = '  +2.43%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.405'
$ws.Range("D28").Style = "Normal"

$ws.Range("E28").Value = '  +0.36%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '106.66'
$ws.Range("D29").Style = "Normal"

$ws.Range("E29").Value = '  -1.48%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.984'
$ws.Range("D30").Style = "Normal"

$ws.Range("E30").Value = '  -0.20%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.07966'
$ws.Range("D31").Style = "Normal"

$ws.Range("E31").Value = '  +2.74%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.713'
$ws.Range("D32").Style = "Normal"

$ws.Range("E32").Value = '  +0.84%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04585'
$ws.Range("D33").Style = "Normal"

$ws.Range("E33").Value = '  +3.54%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.000'
$ws.Range("D34").Style = "Normal"

$ws.Range("E34").Value = '  +0.02%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.620'
$ws.Range("D35").Style = "Normal"

$ws.Range("E35").Value = '  +0.26%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9946'
$ws.Range("D36").Style = "Normal"

$ws.Range("E36").Value = '  +1.93%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.6313'
$ws.Range("D37").Style = "Normal"

$ws.Range("E37").Value = '  +2.17%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.095'
$ws.Range("D38").Style = "Normal"

$ws.Range("E38").Value = '  +9.59%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.9190'
$ws.Range("D39").Style = "Normal"

$ws.Range("E39").Value = '  -1.49%  '

$ws.Range("E40").Value = '  -0.64%  '

$ws.Range("B41").Value = 'Quant'

$ws.Range("C41").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '104.59'
$ws.Range("D41").Style = "Normal"

$ws.Range("E41").Value = '  -7.35%  '

$ws.Range("B42").Value = 'PaxDollar'

$ws.Range("C42").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.007'
$ws.Range("D42").Style = "Normal"

$ws.Range("E42").Value = '  +0.64%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.01504'
$ws.Range("D43").Style = "Normal"

$ws.Range("E43").Value = '  +2.10%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.583'
$ws.Range("D44").Style = "Normal"

$ws.Range("E44").Value = '  +4.82%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.3875'
$ws.Range("D45").Style = "Normal"

$ws.Range("E45").Value = '  +1.93%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '6.961'
$ws.Range("D46").Style = "Normal"

$ws.Range("E46").Value = '  +10.56%  '

$ws.Range("E47").Value = '  +1.54%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.05348'
$ws.Range("D48").Style = "Normal"

$ws.Range("E48").Value = '  +1.23%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '30.93'
$ws.Range("D49").Style = "Normal"

$ws.Range("E49").Value = '  +2.13%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.793'
$ws.Range("D50").Style = "Normal"

$ws.Range("E50").Value = '  +1.33%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.258'
$ws.Range("D51").Style = "Normal"

$ws.Range("E51").Value = '  +3.89%  '
